# Generate Report for Handback
# Updates the localization-status report: the handback is now in sync with
# en-US (status text changes), the "Latest Handback DateTime" values move
# forward to the new generation timestamps, and the per-language "Error
# Detail" messages are cleared now that the handback versions are current.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$newStatus = "Handed back: in sync with en-US"

# --- Overview sheet: zh-cn / de-de status columns ---
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus

# --- zh-cn detail sheet ---
$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("K2").Value = "2016-08-12 19:01:30"
$zhcn.Range("P2").Value = ""

# --- de-de detail sheet ---
$dede.Range("C2").Value = $newStatus
$dede.Range("K2").Value = "2016-08-12 19:01:40"
$dede.Range("P2").Value = ""

# --- Column widths: widen the Status columns to fit the longer text, and
#     shrink the now-empty Error Detail columns back down. The host's
#     ColumnWidth setter snaps to a whole-pixel grid, so we pick the value
#     (the midpoint of the input range that rounds to the desired pixel
#     width) closest to the authored widths of 29.9777047293527 /
#     13.7470528738839.
$statusColWidth = 29.16666666666667
$errorColWidth  = 12.83333333333333

$overview.Columns.Item(5).ColumnWidth = $statusColWidth
$overview.Columns.Item(6).ColumnWidth = $statusColWidth

$zhcn.Columns.Item(3).ColumnWidth  = $statusColWidth
$zhcn.Columns.Item(16).ColumnWidth = $errorColWidth

$dede.Columns.Item(3).ColumnWidth  = $statusColWidth
$dede.Columns.Item(16).ColumnWidth = $errorColWidth

Write-Host "Applied handback report updates"
